$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$helper = $ws.Range("Z1")
$helper.NumberFormat = "@"

function Set-TextValue($cellRef, $val) {
    $helper.Value = $val
    $helper.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
}

Set-TextValue 'D2' '25.851.13'
Set-TextValue 'E2' '  +0.23%  '
Set-TextValue 'D3' '1.630.01'
Set-TextValue 'E3' '  +0.01%  '
Set-TextValue 'E4' '  +0.45%  '
Set-TextValue 'D5' '214.38'
Set-TextValue 'E5' '  +0.16%  '
Set-TextValue 'E6' '  +0.74%  '
Set-TextValue 'E7' '  +0.37%  '
Set-TextValue 'E9' '  -0.05%  '
Set-TextValue 'D10' '19.52'
Set-TextValue 'E10' '  -0.62%  '
Set-TextValue 'D11' '0.0789'
Set-TextValue 'E11' '  -0.08%  '
Set-TextValue 'D12' '1.855.14'
Set-TextValue 'E12' '  +0.00%  '
Set-TextValue 'E13' '  +0.09%  '
Set-TextValue 'D14' '1.611.41'
Set-TextValue 'E14' '  -1.06%  '
Set-TextValue 'D15' '0.543'
Set-TextValue 'E15' '  -1.56%  '
Set-TextValue 'D16' '0.0₃0753'
Set-TextValue 'E16' '  -0.86%  '
Set-TextValue 'D17' '62.48'
Set-TextValue 'E17' '  -0.39%  '
Set-TextValue 'D18' '25.852.72'
Set-TextValue 'E18' '  +0.26%  '
Set-TextValue 'E19' '  +0.39%  '
Set-TextValue 'B20' 'BitcoinCash'
Set-TextValue 'C20' 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 'D20' '193.14'
Set-TextValue 'E20' '  +1.23%  '
Set-TextValue 'B21' 'Uniswap'
Set-TextValue 'C21' 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 'D21' '4.38'
Set-TextValue 'E21' '  -1.03%  '
Set-TextValue 'E22' '  +0.00%  '
Set-TextValue 'E23' '  -0.43%  '
Set-TextValue 'D24' '1.81'
Set-TextValue 'E24' '  -0.08%  '
Set-TextValue 'B25' 'BinanceUSD'
Set-TextValue 'C25' 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue 'D25' '1.00'
Set-TextValue 'E25' '  +0.44%  '
Set-TextValue 'B26' 'Monero'
Set-TextValue 'C26' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D26' '143.00'
Set-TextValue 'E26' '  +0.64%  '
Set-TextValue 'E27' '  +2.81%  '
Set-TextValue 'E28' '  +0.03%  '
Set-TextValue 'D29' '15.41'
Set-TextValue 'E29' '  -0.51%  '
Set-TextValue 'E30' '  -0.03%  '
Set-TextValue 'D31' '0.0497'
Set-TextValue 'E31' '  +0.48%  '
Set-TextValue 'E33' '  +0.02%  '
Set-TextValue 'D34' '1.56'
Set-TextValue 'E34' '  -1.28%  '
Set-TextValue 'E35' '  +1.49%  '
Set-TextValue 'D36' '0.899'
Set-TextValue 'E36' '  -0.51%  '
Set-TextValue 'D37' '1.136.28'
Set-TextValue 'E37' '  -0.36%  '
Set-TextValue 'E38' '  +0.40%  '
Set-TextValue 'D39' '2.47'
Set-TextValue 'E39' '  -1.36%  '
Set-TextValue 'E40' '  +0.29%  '
Set-TextValue 'E41' '  +0.48%  '
Set-TextValue 'B42' 'TrustWalletToken'
Set-TextValue 'C42' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D42' '0.801'
Set-TextValue 'E42' '  -0.02%  '
Set-TextValue 'B43' 'Quant'
Set-TextValue 'C43' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 'D43' '99.03'
Set-TextValue 'E43' '  -1.78%  '
Set-TextValue 'B44' 'FraxShare'
Set-TextValue 'C44' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D44' '5.44'
Set-TextValue 'E44' '  -2.66%  '
Set-TextValue 'D45' '1.764.49'
Set-TextValue 'E45' '  -0.01%  '
Set-TextValue 'B46' 'Aave'
Set-TextValue 'C46' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D46' '56.17'
Set-TextValue 'E46' '  +1.60%  '
Set-TextValue 'B47' 'Cronos'
Set-TextValue 'C47' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D47' '0.0525'
Set-TextValue 'E47' '  +2.64%  '
Set-TextValue 'B48' 'RenderToken'
Set-TextValue 'C48' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D48' '1.45'
Set-TextValue 'E48' '  -1.49%  '
Set-TextValue 'B49' 'Mantle'
Set-TextValue 'C49' 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 'D49' '0.415'
Set-TextValue 'E49' '  -0.36%  '
Set-TextValue 'B50' 'EnergySwap'
Set-TextValue 'C50' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D50' '7.59'
Set-TextValue 'E50' '  +0.84%  '
Set-TextValue 'B51' 'Algorand'
Set-TextValue 'C51' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 'D51' '0.0958'
Set-TextValue 'E51' '  +0.66%  '

$helper.Clear()
$excel.CutCopyMode = $false
